$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "三花智控"
$ws.Range("A3").Value = "山子高科"
$ws.Range("C3").Value = "山子高科"
$ws.Range("A4").Value = "首开股份"
$ws.Range("C4").Value = "上海建工"
$ws.Range("A5").Value = "上海建工"
$ws.Range("B5").Value = "三花智控"
$ws.Range("B6").Value = "首开股份"
$ws.Range("C6").Value = "华胜天成"
$ws.Range("A7").Value = "均胜电子"
$ws.Range("B7").Value = "金发科技"
$ws.Range("C7").Value = "利欧股份"
$ws.Range("A8").Value = "金发科技"
$ws.Range("B8").Value = "工业富联"
$ws.Range("C8").Value = "首开股份"
$ws.Range("A9").Value = "万向钱潮"
$ws.Range("B9").Value = "利欧股份"
$ws.Range("C9").Value = "金发科技"
$ws.Range("A10").Value = "工业富联"
$ws.Range("B10").Value = "露笑科技"
$ws.Range("C10").Value = "万通发展"
$ws.Range("A11").Value = "华胜天成"
$ws.Range("B11").Value = "华胜天成"
$ws.Range("A12").Value = "利欧股份"
$ws.Range("B12").Value = "均胜电子"
$ws.Range("C12").Value = "天赐材料"
$ws.Range("A13").Value = "景兴纸业"
$ws.Range("B13").Value = "巨轮智能"
$ws.Range("C13").Value = "万向钱潮"
$ws.Range("A14").Value = "露笑科技"
$ws.Range("B14").Value = "万向钱潮"
$ws.Range("C14").Value = "先导智能"
$ws.Range("A15").Value = "巨轮智能"
$ws.Range("B15").Value = "东方财富"
$ws.Range("C15").Value = "露笑科技"
$ws.Range("A16").Value = "万通发展"
$ws.Range("B16").Value = "景兴纸业"
$ws.Range("C16").Value = "巨轮智能"
$ws.Range("A17").Value = "东方财富"
$ws.Range("B17").Value = "和而泰"
$ws.Range("C17").Value = "东方财富"
$ws.Range("A18").Value = "凯美特气"
$ws.Range("B18").Value = "先导智能"
$ws.Range("C18").Value = "青山纸业"
$ws.Range("A19").Value = "寒武纪-U"
$ws.Range("B19").Value = "万通发展"
$ws.Range("C19").Value = "拓维信息"
$ws.Range("A20").Value = "中芯国际"
$ws.Range("B20").Value = "青山纸业"
$ws.Range("C20").Value = "凯美特气"
$ws.Range("A21").Value = "先导智能"
$ws.Range("B21").Value = "凯美特气"
$ws.Range("C21").Value = "岩山科技"
